# Updates cryptos list figures (price / 1h volume) to the latest scrape.
# Row 44/45 additionally swap rank position between Mantle and WhiteBITCoin.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "59.251.42"
$ws.Range("E2").Value = "  +1.98%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.586.45"
$ws.Range("E3").Value = "  -0.27%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D5").Value = "523.52"
$ws.Range("E5").Value = "  +0.29%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D6").Value = "139.20"
$ws.Range("E6").Value = "  -3.11%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.17%  "

# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  -0.73%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.597.18"
$ws.Range("E9").Value = "  -0.61%  "

# Row 10: Toncoin
$ws.Range("D10").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D10").Value = "6.51"
$ws.Range("E10").Value = "  +0.10%  "

# Row 11: Dogecoin
$ws.Range("E11").Value = "  -0.28%  "

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D12").Value = "0.331"
$ws.Range("E12").Value = "  -1.89%  "

# Row 13: TRON
$ws.Range("E13").Value = "  +2.77%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.047.02"
$ws.Range("E14").Value = "  -0.17%  "

# Row 15: WrappedBTC
$ws.Range("D15").Value = "59.172.95"
$ws.Range("E15").Value = "  +1.86%  "

# Row 16: Avalanche
$ws.Range("D16").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D16").Value = "20.46"
$ws.Range("E16").Value = "  +0.11%  "

# Row 17: ShibaInu
$ws.Range("E17").Value = "  -0.74%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "2.581.27"
$ws.Range("E18").Value = "  +1.29%  "

# Row 19: BitcoinCash
$ws.Range("D19").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D19").Value = "341.77"
$ws.Range("E19").Value = "  +0.80%  "

# Row 20: Polkadot
$ws.Range("E20").Value = "  -1.16%  "

# Row 21: Chainlink
$ws.Range("E21").Value = "  -1.92%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D22").Value = "6.44"
$ws.Range("E22").Value = "  +0.48%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.13%  "

# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D24").Value = "66.90"
$ws.Range("E24").Value = "  +2.37%  "

# Row 25: Kaspa
$ws.Range("E25").Value = "  +0.88%  "

# Row 26: Polygon
$ws.Range("D26").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D26").Value = "0.405"
$ws.Range("E26").Value = "  +0.65%  "

# Row 27: Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.04%  "

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D28").Value = "7.06"
$ws.Range("E28").Value = "  +0.72%  "

# Row 30: PEPE
$ws.Range("D30").Value = "0.0₃0723"

# Row 31: Aptos
$ws.Range("E31").Value = "  -4.96%  "

# Row 32: PancakeSwap
$ws.Range("E32").Value = "  +0.13%  "

# Row 33: EthereumClassic
$ws.Range("D33").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D33").Value = "18.71"
$ws.Range("E33").Value = "  -0.30%  "

# Row 34: Monero
$ws.Range("D34").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D34").Value = "149.19"
$ws.Range("E34").Value = "  -0.37%  "

# Row 35: NEARProtocol
$ws.Range("D35").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D35").Value = "3.97"
$ws.Range("E35").Value = "  -1.33%  "

# Row 36: ImmutableX
$ws.Range("E36").Value = "  -1.59%  "

# Row 37: OKB
$ws.Range("D37").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D37").Value = "36.79"
$ws.Range("E37").Value = "  +2.28%  "

# Row 38: Stacks
$ws.Range("D38").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D38").Value = "1.47"
$ws.Range("E38").Value = "  +1.63%  "

# Row 39: Fetch.AI
$ws.Range("D39").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  -4.30%  "

# Row 40: SuiNetwork
$ws.Range("D40").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").Value = "  -6.54%  "

# Row 41: Filecoin
$ws.Range("E41").Value = "  -0.74%  "

# Row 42: FirstDigitalUSD
$ws.Range("E42").Value = "  +0.26%  "

# Row 43: Bittensor
$ws.Range("D43").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D43").Value = "270.73"
$ws.Range("E43").Value = "  -0.50%  "

# Row 44: Mantle
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D44").Value = "10.80"
$ws.Range("E44").Value = "  +1.19%  "

# Row 45: WhiteBITCoin
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D45").Value = "0.596"
$ws.Range("E45").Value = "  -0.62%  "

# Row 46: Stellar
$ws.Range("E46").Value = "  -0.67%  "

# Row 47: Hedera
$ws.Range("D47").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D47").Value = "0.0513"
$ws.Range("E47").Value = "  -1.68%  "

# Row 48: EnergySwap
$ws.Range("D48").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D48").Value = "18.39"
$ws.Range("E48").Value = "  -2.11%  "

# Row 49: Maker
$ws.Range("D49").Value = "1.969.28"
$ws.Range("E49").Value = "  -0.30%  "

# Row 50: VeChain
$ws.Range("E50").Value = "  +0.19%  "

# Row 51: InjectiveProtocol
$ws.Range("D51").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D51").Value = "18.05"
$ws.Range("E51").Value = "  -4.37%  "
